$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the gate closure trigger values based on feedback from CPRA
$ws.Range("C1").Value = 4
$ws.Range("C9").Value = 2
$ws.Range("C11").Value = 2.5

# Update the active cell/selection on the sheet view to D1
$ws.Range("D1").Select()
